$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 17957.105
$ws.Range("J62").Value = 14888.556
$ws.Range("L62").Value = 14888.556
$ws.Range("N62").Value = -16136.556

# Row 65
$ws.Range("H65").Value = 17957.105
$ws.Range("J65").Value = 14888.556
$ws.Range("L65").Value = 74442.78
$ws.Range("N65").Value = -80682.78

# Row 138
$ws.Range("H138").Value = 4411
$ws.Range("I138").Value = 1668.8572
$ws.Range("J138").Value = 4903.1797
$ws.Range("K138").Value = 5006.571599999999
$ws.Range("L138").Value = 14709.5391
$ws.Range("M138").Value = 133.4284000000007
$ws.Range("N138").Value = -24989.5391

# Row 141
$ws.Range("H141").Value = 1533.4166
$ws.Range("I141").Value = 1533.4166
$ws.Range("K141").Value = 4600.2498
$ws.Range("M141").Value = 579.7502000000004

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 886.0897
$ws.Range("I32").Value = 805.2089999999999
$ws.Range("K32").Value = 805.2089999999999
$ws.Range("M32").Value = -518.2089999999999

# Row 33
$ws.Range("H33").Value = 9474.75
$ws.Range("I33").Value = 9999.5
$ws.Range("J33").Value = 8950
$ws.Range("K33").Value = 9999.5
$ws.Range("L33").Value = 8950
$ws.Range("M33").Value = -9670.5
$ws.Range("N33").Value = -9608

# Row 76
$ws.Range("H76").Value = 47599.8
$ws.Range("J76").Value = 47599.8
$ws.Range("L76").Value = 47599.8
$ws.Range("N76").Value = -48275.8

# Row 79
$ws.Range("H79").Value = 47599.8
$ws.Range("J79").Value = 47599.8
$ws.Range("L79").Value = 47599.8
$ws.Range("N79").Value = -49939.8

# Row 88
$ws.Range("H88").Value = 1708.4166
$ws.Range("I88").Value = 1701.25
$ws.Range("J88").Value = 1712
$ws.Range("K88").Value = 1701.25
$ws.Range("L88").Value = 1712
$ws.Range("M88").Value = -1295.25
$ws.Range("N88").Value = -2524

# Row 91
$ws.Range("H91").Value = 1708.4166
$ws.Range("I91").Value = 1701.25
$ws.Range("J91").Value = 1712
$ws.Range("K91").Value = 1701.25
$ws.Range("L91").Value = 1712
$ws.Range("M91").Value = -297.25
$ws.Range("N91").Value = -4520

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1818.6111
$ws.Range("I86").Value = 1719.8667
$ws.Range("K86").Value = 1719.8667
$ws.Range("M86").Value = -596.8667

# Row 89
$ws.Range("H89").Value = 1818.6111
$ws.Range("I89").Value = 1719.8667
$ws.Range("K89").Value = 8599.333500000001
$ws.Range("M89").Value = -2983.333500000001

# Row 97
$ws.Range("H97").Value = 3459.182
$ws.Range("I97").Value = 2578.7144
$ws.Range("K97").Value = 2578.7144
$ws.Range("M97").Value = -1587.7144

# Row 99
$ws.Range("H99").Value = 87785
$ws.Range("I99").Value = 87785
$ws.Range("K99").Value = 87785
$ws.Range("M99").Value = -86287

# Row 134
$ws.Range("H134").Value = 5513.6665
$ws.Range("I134").Value = 5445.0356
$ws.Range("J134").Value = 6474.5
$ws.Range("K134").Value = 16335.1068
$ws.Range("L134").Value = 19423.5
$ws.Range("M134").Value = -13800.1068
$ws.Range("N134").Value = -24493.5

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 1674483.4
$ws.Range("I4").Value = 2507000
$ws.Range("K4").Value = 2507000
$ws.Range("M4").Value = -2506888

# Row 22
$ws.Range("H22").Value = 717.6539
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -1000

# Row 35
$ws.Range("H35").Value = 10615.4375
$ws.Range("I35").Value = 9680.615
$ws.Range("J35").Value = 14666.333
$ws.Range("K35").Value = 9680.615
$ws.Range("L35").Value = 14666.333
$ws.Range("M35").Value = -9386.615
$ws.Range("N35").Value = -15254.333

# Row 86
$ws.Range("H86").Value = 6176.857
$ws.Range("I86").Value = 4822.75
$ws.Range("K86").Value = 4822.75
$ws.Range("M86").Value = -3699.75

# Row 89
$ws.Range("H89").Value = 6176.857
$ws.Range("I89").Value = 4822.75
$ws.Range("K89").Value = 24113.75
$ws.Range("M89").Value = -18497.75

# Row 99
$ws.Range("H99").Value = 5124.625
$ws.Range("I99").Value = 5249.8335
$ws.Range("J99").Value = 4749
$ws.Range("K99").Value = 5249.8335
$ws.Range("L99").Value = 4749
$ws.Range("M99").Value = -3751.8335
$ws.Range("N99").Value = -7745

# Row 126
$ws.Range("H126").Value = 5124.625
$ws.Range("I126").Value = 5249.8335
$ws.Range("J126").Value = 4749
$ws.Range("K126").Value = 15749.5005
$ws.Range("L126").Value = 14247
$ws.Range("M126").Value = -13279.5005
$ws.Range("N126").Value = -19187

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 1594.5
$ws.Range("I23").Value = 831.4
$ws.Range("J23").Value = 2866.3333
$ws.Range("K23").Value = 2494.2
$ws.Range("L23").Value = 8598.999899999999
$ws.Range("M23").Value = -2259.2
$ws.Range("N23").Value = -9068.999899999999

# Row 113
$ws.Range("H113").Value = 600
$ws.Range("I113").Value = 300
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270

# Row 114
$ws.Range("H114").Value = 3564.9546
$ws.Range("J114").Value = 4398.625
$ws.Range("L114").Value = 13195.875
$ws.Range("N114").Value = -19703.875

# Row 126
$ws.Range("H126").Value = 1974.75
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -15880

# Row 141
$ws.Range("H141").Value = 5550.75
$ws.Range("I141").Value = 2656.125
$ws.Range("K141").Value = 7968.375
$ws.Range("M141").Value = -2788.375

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 297.85715
$ws.Range("I2").Value = 343.33334
$ws.Range("K2").Value = 343.33334
$ws.Range("M2").Value = -230.33334

# Row 36
$ws.Range("H36").Value = 3971.2856
$ws.Range("J36").Value = 4239.8
$ws.Range("L36").Value = 4239.8
$ws.Range("N36").Value = -5209.8

# Row 80
$ws.Range("H80").Value = 84086.25
$ws.Range("I80").Value = 250850
$ws.Range("J80").Value = 28498.334
$ws.Range("K80").Value = 250850
$ws.Range("L80").Value = 28498.334
$ws.Range("M80").Value = -249852
$ws.Range("N80").Value = -30494.334

# Row 83
$ws.Range("H83").Value = 84086.25
$ws.Range("I83").Value = 250850
$ws.Range("J83").Value = 28498.334
$ws.Range("K83").Value = 1254250
$ws.Range("L83").Value = 142491.67
$ws.Range("M83").Value = -1249258
$ws.Range("N83").Value = -152475.67

# Row 134
$ws.Range("H134").Value = 49997.8
$ws.Range("J134").Value = 49997.8
$ws.Range("L134").Value = 149993.4
$ws.Range("N134").Value = -155063.4

# Row 136
$ws.Range("H136").Value = 49999
$ws.Range("J136").Value = 49999
$ws.Range("L136").Value = 149997
$ws.Range("N136").Value = -155097

$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Range("H74").Value = 76173.60000000001
$ws.Range("J74").Value = 76173.60000000001
$ws.Range("L74").Value = 76173.60000000001
$ws.Range("N74").Value = -78169.60000000001

# Row 77
$ws.Range("H77").Value = 76173.60000000001
$ws.Range("J77").Value = 76173.60000000001
$ws.Range("L77").Value = 228520.8
$ws.Range("N77").Value = -238504.8

# Row 132
$ws.Range("H132").Value = 3706.48
$ws.Range("I132").Value = 3731.6667
$ws.Range("K132").Value = 11195.0001
$ws.Range("M132").Value = -8665.000100000001

# Row 136
$ws.Range("H136").Value = 35464.867
$ws.Range("I136").Value = 3559.8572
$ws.Range("K136").Value = 10679.5716
$ws.Range("M136").Value = -8129.571599999999

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 27605.467
$ws.Range("I62").Value = 3987.25
$ws.Range("J62").Value = 36193.91
$ws.Range("K62").Value = 3987.25
$ws.Range("L62").Value = 36193.91
$ws.Range("M62").Value = -3363.25
$ws.Range("N62").Value = -37441.91

# Row 65
$ws.Range("H65").Value = 27605.467
$ws.Range("I65").Value = 3987.25
$ws.Range("J65").Value = 36193.91
$ws.Range("K65").Value = 19936.25
$ws.Range("L65").Value = 180969.55
$ws.Range("M65").Value = -16816.25
$ws.Range("N65").Value = -187209.55

# Row 136
$ws.Range("H136").Value = 8148.6445
$ws.Range("I136").Value = 9651.267
$ws.Range("K136").Value = 28953.801
$ws.Range("M136").Value = -26403.801

# Row 137
$ws.Range("H137").Value = 74499.5
$ws.Range("J137").Value = 74499.5
$ws.Range("L137").Value = 74499.5
$ws.Range("N137").Value = -84699.5

# Row 138
$ws.Range("H138").Value = 76999.5
$ws.Range("J138").Value = 76999.5
$ws.Range("L138").Value = 76999.5
$ws.Range("N138").Value = -87279.5
